# "fixing contd.. in b slides"
#
# Slides 1-3 each carry a small "contd.." caption textbox
# (Google Shape;131;p25) in the bottom-right corner. Re-touching that
# textbox's paragraph/bullet formatting (as happened when the author
# re-selected it in PowerPoint) causes PowerPoint to write out the
# paragraph properties explicitly instead of relying on inherited
# defaults: explicit left margin/indent/level/alignment/direction,
# explicit (zero) space-before/after, and an explicit "no bullet"
# (Arial bullet font, bullet turned off). Slide 4's matching shape has
# no "contd.." run, so it is left untouched.

$p = $ppt.ActivePresentation

for ($slideIndex = 1; $slideIndex -le 3; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(3)              # Google Shape;131;p25 ("contd..")

    $textRange = $shape.TextFrame.TextRange
    $paraFormat = $textRange.ParagraphFormat

    # marL="0" / indent="0" via the paragraph's ruler level.
    $rulerLevel = $shape.TextFrame.Ruler.Levels.Item(1)
    $rulerLevel.LeftMargin = 0
    $rulerLevel.FirstMargin = 0

    # lvl="0" (IndentLevel is 1-based; level 1 -> lvl="0").
    $textRange.IndentLevel = 1

    # algn="l"
    $paraFormat.Alignment = 1

    # rtl="0"
    $paraFormat.TextDirection = 1

    # spcBef/spcAft val="0" (DrawingML-level paragraph format).
    $tr2 = $shape.TextFrame2.TextRange
    $tr2.ParagraphFormat.SpaceBefore = 0
    $tr2.ParagraphFormat.SpaceAfter = 0

    # buFont typeface="Arial" + buNone
    $paraFormat.Bullet.Font.Name = "Arial"
    $paraFormat.Bullet.Visible = $false
}
